$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D column (last "未充电截止时间" recalculation timestamp) for rows 2-18
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 4).Value = 46007.305775462963
}

# Update rows 19-59 with refreshed data (A: station name, B: terminal name, C: last charge end time, D: cutoff time)
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "603号直流"
$ws.Cells.Item(19, 3).Value = 46003.262175925927
$ws.Cells.Item(19, 4).Value = 46007.305775462963
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(20, 2).Value = "101号直流"
$ws.Cells.Item(20, 3).Value = 46004.217581018522
$ws.Cells.Item(20, 4).Value = 46007.305775462963
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "505号直流"
$ws.Cells.Item(21, 3).Value = 46004.540798611109
$ws.Cells.Item(21, 4).Value = 46007.305775462963
$ws.Cells.Item(22, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(22, 2).Value = "303号直流"
$ws.Cells.Item(22, 3).Value = 46004.548645833333
$ws.Cells.Item(22, 4).Value = 46007.305775462963
$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23, 2).Value = "904号直流"
$ws.Cells.Item(23, 3).Value = 46005.57613425926
$ws.Cells.Item(23, 4).Value = 46007.305775462963
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "401号直流"
$ws.Cells.Item(24, 3).Value = 46005.688252314816
$ws.Cells.Item(24, 4).Value = 46007.305775462963
$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(25, 2).Value = "304号直流"
$ws.Cells.Item(25, 3).Value = 46006.083252314813
$ws.Cells.Item(25, 4).Value = 46007.305775462963
$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(26, 2).Value = "007A号直流"
$ws.Cells.Item(26, 3).Value = 46006.123229166667
$ws.Cells.Item(26, 4).Value = 46007.305775462963
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(27, 2).Value = "905号直流"
$ws.Cells.Item(27, 3).Value = 46006.212766203702
$ws.Cells.Item(27, 4).Value = 46007.305775462963
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(28, 2).Value = "802号直流"
$ws.Cells.Item(28, 3).Value = 46006.252256944441
$ws.Cells.Item(28, 4).Value = 46007.305775462963
$ws.Cells.Item(29, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(29, 2).Value = "208号直流"
$ws.Cells.Item(29, 3).Value = 46006.461678240739
$ws.Cells.Item(29, 4).Value = 46007.305775462963
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30, 2).Value = "103号直流"
$ws.Cells.Item(30, 3).Value = 46006.528240740743
$ws.Cells.Item(30, 4).Value = 46007.305775462963
$ws.Cells.Item(31, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "604号直流"
$ws.Cells.Item(31, 3).Value = 46006.5391087963
$ws.Cells.Item(31, 4).Value = 46007.305775462963
$ws.Cells.Item(32, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(32, 2).Value = "005A号直流"
$ws.Cells.Item(32, 3).Value = 46006.553807870368
$ws.Cells.Item(32, 4).Value = 46007.305775462963
$ws.Cells.Item(33, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value = "301号直流"
$ws.Cells.Item(33, 3).Value = 46006.555787037039
$ws.Cells.Item(33, 4).Value = 46007.305775462963
$ws.Cells.Item(34, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(34, 2).Value = "204号直流"
$ws.Cells.Item(34, 3).Value = 46006.557673611111
$ws.Cells.Item(34, 4).Value = 46007.305775462963
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value = "001A号直流"
$ws.Cells.Item(35, 3).Value = 46006.559583333335
$ws.Cells.Item(35, 4).Value = 46007.305775462963
$ws.Cells.Item(36, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(36, 2).Value = "003B号直流"
$ws.Cells.Item(36, 3).Value = 46006.562893518516
$ws.Cells.Item(36, 4).Value = 46007.305775462963
$ws.Cells.Item(37, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(37, 2).Value = "702号直流"
$ws.Cells.Item(37, 3).Value = 46006.565613425926
$ws.Cells.Item(37, 4).Value = 46007.305775462963
$ws.Cells.Item(38, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(38, 2).Value = "110号直流"
$ws.Cells.Item(38, 3).Value = 46006.565717592595
$ws.Cells.Item(38, 4).Value = 46007.305775462963
$ws.Cells.Item(39, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(39, 2).Value = "503号直流"
$ws.Cells.Item(39, 3).Value = 46006.573784722219
$ws.Cells.Item(39, 4).Value = 46007.305775462963
$ws.Cells.Item(40, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(40, 2).Value = "107号直流"
$ws.Cells.Item(40, 3).Value = 46006.575752314813
$ws.Cells.Item(40, 4).Value = 46007.305775462963
$ws.Cells.Item(41, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(41, 2).Value = "011A号直流"
$ws.Cells.Item(41, 3).Value = 46006.581087962964
$ws.Cells.Item(41, 4).Value = 46007.305775462963
$ws.Cells.Item(42, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(42, 2).Value = "B01号直流"
$ws.Cells.Item(42, 3).Value = 46006.581782407404
$ws.Cells.Item(42, 4).Value = 46007.305775462963
$ws.Cells.Item(43, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(43, 2).Value = "201号直流"
$ws.Cells.Item(43, 3).Value = 46006.584479166668
$ws.Cells.Item(43, 4).Value = 46007.305775462963
$ws.Cells.Item(44, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(44, 2).Value = "405号直流"
$ws.Cells.Item(44, 3).Value = 46006.589328703703
$ws.Cells.Item(44, 4).Value = 46007.305775462963
$ws.Cells.Item(45, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(45, 2).Value = "903号直流"
$ws.Cells.Item(45, 3).Value = 46006.601030092592
$ws.Cells.Item(45, 4).Value = 46007.305775462963
$ws.Cells.Item(46, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(46, 2).Value = "805号直流"
$ws.Cells.Item(46, 3).Value = 46006.610844907409
$ws.Cells.Item(46, 4).Value = 46007.305775462963
$ws.Cells.Item(47, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(47, 2).Value = "002B号直流"
$ws.Cells.Item(47, 3).Value = 46006.625543981485
$ws.Cells.Item(47, 4).Value = 46007.305775462963
$ws.Cells.Item(48, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(48, 2).Value = "109号直流"
$ws.Cells.Item(48, 3).Value = 46006.629062499997
$ws.Cells.Item(48, 4).Value = 46007.305775462963
$ws.Cells.Item(49, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(49, 2).Value = "212号直流"
$ws.Cells.Item(49, 3).Value = 46006.629502314812
$ws.Cells.Item(49, 4).Value = 46007.305775462963
$ws.Cells.Item(50, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(50, 2).Value = "803号直流"
$ws.Cells.Item(50, 3).Value = 46006.640960648147
$ws.Cells.Item(50, 4).Value = 46007.305775462963
$ws.Cells.Item(51, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(51, 2).Value = "305号直流"
$ws.Cells.Item(51, 3).Value = 46006.646874999999
$ws.Cells.Item(51, 4).Value = 46007.305775462963
$ws.Cells.Item(52, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(52, 2).Value = "503号直流"
$ws.Cells.Item(52, 3).Value = 46006.657731481479
$ws.Cells.Item(52, 4).Value = 46007.305775462963
$ws.Cells.Item(53, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(53, 2).Value = "101号直流"
$ws.Cells.Item(53, 3).Value = 46006.687800925924
$ws.Cells.Item(53, 4).Value = 46007.305775462963
$ws.Cells.Item(54, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(54, 2).Value = "403号直流"
$ws.Cells.Item(54, 3).Value = 46006.68787037037
$ws.Cells.Item(54, 4).Value = 46007.305775462963
$ws.Cells.Item(55, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(55, 2).Value = "004A号直流"
$ws.Cells.Item(55, 3).Value = 46006.692013888889
$ws.Cells.Item(55, 4).Value = 46007.305775462963
$ws.Cells.Item(56, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(56, 2).Value = "008B号直流"
$ws.Cells.Item(56, 3).Value = 46006.70590277778
$ws.Cells.Item(56, 4).Value = 46007.305775462963
$ws.Cells.Item(57, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(57, 2).Value = "108号直流"
$ws.Cells.Item(57, 3).Value = 46006.727013888885
$ws.Cells.Item(57, 4).Value = 46007.305775462963
$ws.Cells.Item(58, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(58, 2).Value = "203号直流"
$ws.Cells.Item(58, 3).Value = 46006.738391203704
$ws.Cells.Item(58, 4).Value = 46007.305775462963
$ws.Cells.Item(59, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(59, 2).Value = "002A号直流"
$ws.Cells.Item(59, 3).Value = 46006.775914351849
$ws.Cells.Item(59, 4).Value = 46007.305775462963

# Update selected cell in sheet view
$ws.Range("G18").Select()
